$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 38039
$ws1.Cells.Item(4, 6).Value = 644
$ws1.Cells.Item(6, 6).Value = 493
$ws1.Cells.Item(10, 6).Value = 112
$ws1.Cells.Item(11, 6).Value = 767
$ws1.Cells.Item(12, 6).Value = 601
$ws1.Cells.Item(13, 6).Value = 93
$ws1.Cells.Item(16, 6).Value = 700
$ws1.Cells.Item(17, 6).Value = 197
$ws1.Cells.Item(18, 6).Value = 500
$ws1.Cells.Item(20, 6).Value = 1200
$ws1.Cells.Item(22, 6).Value = 893
$ws1.Cells.Item(23, 6).Value = 2618
$ws1.Cells.Item(24, 6).Value = 1106
$ws1.Cells.Item(26, 6).Value = 132
$ws1.Cells.Item(29, 6).Value = 858
$ws1.Cells.Item(30, 6).Value = 81
$ws1.Cells.Item(31, 6).Value = 1199

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 467
$ws2.Cells.Item(4, 6).Value = 339
$ws2.Cells.Item(7, 6).Value = 57
$ws2.Cells.Item(9, 6).Value = 145
$ws2.Cells.Item(3, 7).Value = 233

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 683

$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 683
$ws4.Cells.Item(3, 6).Value = 38039
$ws4.Cells.Item(5, 6).Value = 644
$ws4.Cells.Item(7, 6).Value = 493
$ws4.Cells.Item(11, 6).Value = 467
$ws4.Cells.Item(12, 6).Value = 467
$ws4.Cells.Item(13, 6).Value = 339
$ws4.Cells.Item(17, 6).Value = 112
$ws4.Cells.Item(18, 6).Value = 767
$ws4.Cells.Item(19, 6).Value = 601
$ws4.Cells.Item(20, 6).Value = 93
$ws4.Cells.Item(21, 6).Value = 57
$ws4.Cells.Item(24, 6).Value = 145
$ws4.Cells.Item(26, 6).Value = 47
$ws4.Cells.Item(28, 6).Value = 700
$ws4.Cells.Item(29, 6).Value = 197
$ws4.Cells.Item(30, 6).Value = 500
$ws4.Cells.Item(32, 6).Value = 1200
$ws4.Cells.Item(34, 6).Value = 893
$ws4.Cells.Item(35, 6).Value = 2618
$ws4.Cells.Item(36, 6).Value = 1106
$ws4.Cells.Item(38, 6).Value = 132
$ws4.Cells.Item(42, 6).Value = 858
$ws4.Cells.Item(43, 6).Value = 81
$ws4.Cells.Item(44, 6).Value = 1200
$ws4.Cells.Item(11, 7).Value = 233
$ws4.Cells.Item(12, 7).Value = 233